# "End of play Saturday"
# The old Question 6 ("Client" / "client" / "blank") table occupied rows 23-27
# and has been removed; everything below it (Question 11 onward) shifts up by
# five rows. Deleting the entire rows (rather than just clearing their
# contents) shifts all the following rows up and keeps every formula/shared
# string reference consistent, exactly like the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 5 rows that made up the old "Client" question block (rows 23-27).
$ws.Rows("23:27").Delete() | Out-Null

# Reflect where the user ended up working afterwards: B23 (the question-11
# "radio" type cell) selected, scrolled down so row 15 is at the top of the
# view.
$ws.Range("B23").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
